# The paper had a trailing run of blank paragraphs (an accidental blank
# "page 5") sitting between the end of the body text and the References
# heading. Collapse that run down to a single blank paragraph so the
# References section moves up right after the body text.

$d = $word.ActiveDocument

# Locate the "References" heading paragraph.
$refIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "References") {
        $refIndex = $i
        break
    }
}

if ($refIndex -gt 1) {
    # Walk backwards from the paragraph just above "References",
    # collecting the contiguous run of paragraphs that contain nothing
    # but whitespace (spaces/tabs) - i.e. the blank "page".
    $lastBlank = $refIndex - 1
    $firstBlank = $lastBlank
    while ($firstBlank -ge 1) {
        $txt = $d.Paragraphs.Item($firstBlank).Range.Text
        $stripped = $txt -replace "[\s\t\r\a]", ""
        if ($stripped.Length -gt 0) {
            break
        }
        $firstBlank = $firstBlank - 1
    }
    $firstBlank = $firstBlank + 1

    # Keep the first blank paragraph as a single spacer, delete the rest
    # of the blank run so References follows right after it.
    if ($lastBlank -gt $firstBlank) {
        $startRange = $d.Paragraphs.Item($firstBlank + 1).Range.Start
        $endRange = $d.Paragraphs.Item($lastBlank).Range.End
        $rng = $d.Range($startRange, $endRange)
        $rng.Delete()
    }
}
